$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 370 ("Fruta / hortaliza, semanal").
# This pushes the previously existing rows 370-461 down by one row (to 371-462),
# and the dimension grows from A1:R461 to A1:R462.
$ws.Rows.Item(370).Insert()

# Populate the newly inserted row 370 with the new record's data.
$ws.Cells.Item(370, 1).Value = 8
$ws.Cells.Item(370, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(370, 3).Value = 'Coquimbo'
$ws.Cells.Item(370, 4).Value = 44964
$ws.Cells.Item(370, 5).Value = 4
$ws.Cells.Item(370, 6).Value = 100112032
$ws.Cells.Item(370, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(370, 8).Value = 'Sin especificar'
$ws.Cells.Item(370, 9).Value = 'Primera'
$ws.Cells.Item(370, 10).Value = 360
$ws.Cells.Item(370, 11).Value = 9000
$ws.Cells.Item(370, 12).Value = 10000
$ws.Cells.Item(370, 13).Value = 9500
$ws.Cells.Item(370, 14).Value = '$/caja 70 unidades'
$ws.Cells.Item(370, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(370, 16).Value = 136
$ws.Cells.Item(370, 17).Value = 70
$ws.Cells.Item(370, 18).Value = 'Hortaliza'
